# Weekly price update: insert a new data row at row 110 (pushing the
# existing rows 110-242 down to 111-243) and populate the new row with
# the latest week's price record for "Macroferia Regional de Talca" /
# Zanahoria.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 110; this shifts all rows
# 110..242 down by one (to 111..243) and keeps their values/formatting.
$ws.Rows.Item(110).Insert()

# Populate the newly inserted row 110 with the new weekly record.
$ws.Range("A110").Value = 5
$ws.Range("B110").Value = "Macroferia Regional de Talca"
$ws.Range("C110").Value = "Maule"
$ws.Range("D110").Value = 44546
$ws.Range("E110").Value = 7
$ws.Range("F110").Value = 100114013
$ws.Range("G110").Value = "Zanahoria"
$ws.Range("H110").Value = "Sin especificar"
$ws.Range("I110").Value = "Primera"
$ws.Range("J110").Value = 400
$ws.Range("K110").Value = 7000
$ws.Range("L110").Value = 7000
$ws.Range("M110").Value = 7000
$ws.Range("N110").Value = "$/saco 20 kilos"
$ws.Range("O110").Value = "Provincia del Elquí"
$ws.Range("P110").Value = 350
$ws.Range("Q110").Value = 20
$ws.Range("R110").Value = "Hortaliza"
